# Update "想去人数" (column F) figures across all sheets to match the
# newly generated data snapshot (gh-pages output regenerated at 456a3b4).
#
# The workbook has 4 sheets:
#   1 - 展览      (Exhibitions)
#   2 - 演出      (Performances)
#   3 - 本地生活  (Local life)
#   4 - 全部类型  (All types - union of sheets 1-3)
# Only column F values change; everything else is untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 781
$ws.Range("F4").Value = 781
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 224
$ws.Range("F7").Value = 35
$ws.Range("F8").Value = 1139
$ws.Range("F9").Value = 895
$ws.Range("F13").Value = 1411
$ws.Range("F14").Value = 55
$ws.Range("F15").Value = 135
$ws.Range("F16").Value = 1584
$ws.Range("F18").Value = 600
$ws.Range("F22").Value = 1077
$ws.Range("F23").Value = 1507
$ws.Range("F24").Value = 745
$ws.Range("F25").Value = 601
$ws.Range("F26").Value = 483
$ws.Range("F27").Value = 467
$ws.Range("F30").Value = 1147
$ws.Range("F31").Value = 295
$ws.Range("F32").Value = 2401
$ws.Range("F33").Value = 274
$ws.Range("F34").Value = 1340
$ws.Range("F35").Value = 456
$ws.Range("F36").Value = 65
$ws.Range("F37").Value = 3925

# --- Sheet 2: 演出 -----------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("F14").Value = 4133
$ws.Range("F18").Value = 43
$ws.Range("F20").Value = 250
$ws.Range("F21").Value = 256
$ws.Range("F25").Value = 43
$ws.Range("F28").Value = 1709

# --- Sheet 3: 本地生活 --------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1264
$ws.Range("F5").Value = 1659
$ws.Range("F8").Value = 986

# --- Sheet 4: 全部类型 (union of the sheets above) ----------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1264
$ws.Range("F4").Value = 1659
$ws.Range("F6").Value = 986
$ws.Range("F8").Value = 781
$ws.Range("F9").Value = 781
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 224
$ws.Range("F12").Value = 35
$ws.Range("F13").Value = 1139
$ws.Range("F14").Value = 895
$ws.Range("F22").Value = 1411
$ws.Range("F23").Value = 55
$ws.Range("F24").Value = 135
$ws.Range("F25").Value = 1584
$ws.Range("F27").Value = 600
$ws.Range("F30").Value = 1077
$ws.Range("F31").Value = 1507
$ws.Range("F33").Value = 745
$ws.Range("F34").Value = 601
$ws.Range("F35").Value = 483
$ws.Range("F36").Value = 467
$ws.Range("F38").Value = 43
$ws.Range("F39").Value = 250
$ws.Range("F40").Value = 256
$ws.Range("F42").Value = 1147
$ws.Range("F43").Value = 295
$ws.Range("F44").Value = 2401
$ws.Range("F46").Value = 43
$ws.Range("F47").Value = 1709
$ws.Range("F48").Value = 1709
$ws.Range("F49").Value = 1340
$ws.Range("F50").Value = 456
$ws.Range("F51").Value = 3925

$wb.Save()
